# Update item_set (B) and sup (C) columns on Sheet1 of the apriori_result workbook.
# B9:B53 shift to the "next" item_set code (one new trailing id appended at B53),
# and every sup value in C2:C57 is refreshed to the latest count.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- item_set (column B), rows 9-53: values shift down by one row, plus a new id at B53 ---
$itemSetUpdates = @(
    @{Row = 9; Value = "21034"},
    @{Row = 10; Value = "21080"},
    @{Row = 11; Value = "21175"},
    @{Row = 12; Value = "21181"},
    @{Row = 13; Value = "21212"},
    @{Row = 14; Value = "21231"},
    @{Row = 15; Value = "21232"},
    @{Row = 16; Value = "21621"},
    @{Row = 17; Value = "21733"},
    @{Row = 18; Value = "21754"},
    @{Row = 19; Value = "21755"},
    @{Row = 20; Value = "21790"},
    @{Row = 21; Value = "21843"},
    @{Row = 22; Value = "21931"},
    @{Row = 23; Value = "21977"},
    @{Row = 24; Value = "22086"},
    @{Row = 25; Value = "22090"},
    @{Row = 26; Value = "22111"},
    @{Row = 27; Value = "22112"},
    @{Row = 28; Value = "22114"},
    @{Row = 29; Value = "22138"},
    @{Row = 30; Value = "22139"},
    @{Row = 31; Value = "22149"},
    @{Row = 32; Value = "22178"},
    @{Row = 33; Value = "22197"},
    @{Row = 34; Value = "22355"},
    @{Row = 35; Value = "22382"},
    @{Row = 36; Value = "22383"},
    @{Row = 37; Value = "22384"},
    @{Row = 38; Value = "22386"},
    @{Row = 39; Value = "22411"},
    @{Row = 40; Value = "22423"},
    @{Row = 41; Value = "22457"},
    @{Row = 42; Value = "22469"},
    @{Row = 43; Value = "22470"},
    @{Row = 44; Value = "47566"},
    @{Row = 45; Value = "48138"},
    @{Row = 46; Value = "82482"},
    @{Row = 47; Value = "82494L"},
    @{Row = 48; Value = "84836"},
    @{Row = 49; Value = "84879"},
    @{Row = 50; Value = "84946"},
    @{Row = 51; Value = "84970S"},
    @{Row = 52; Value = "84991"},
    @{Row = 53; Value = "84992"}
)

foreach ($u in $itemSetUpdates) {
    $cell = $ws.Cells.Item($u.Row, 2)
    # Force text storage so purely-numeric codes (e.g. "21034") are not coerced to numbers,
    # matching the source data where item_set is always a text code.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.ClearFormats()
}

# --- sup (column C), rows 2-57: refreshed support counts ---
$supUpdates = @(
    @{Row = 2; Value = 967},
    @{Row = 3; Value = 900},
    @{Row = 4; Value = 1437},
    @{Row = 5; Value = 847},
    @{Row = 6; Value = 1094},
    @{Row = 7; Value = 878},
    @{Row = 8; Value = 1106},
    @{Row = 9; Value = 1072},
    @{Row = 10; Value = 1078},
    @{Row = 11; Value = 858},
    @{Row = 12; Value = 980},
    @{Row = 13; Value = 1638},
    @{Row = 14; Value = 997},
    @{Row = 15; Value = 1655},
    @{Row = 16; Value = 823},
    @{Row = 17; Value = 1033},
    @{Row = 18; Value = 1302},
    @{Row = 19; Value = 1082},
    @{Row = 20; Value = 959},
    @{Row = 21; Value = 1007},
    @{Row = 22; Value = 1148},
    @{Row = 23; Value = 1046},
    @{Row = 24; Value = 941},
    @{Row = 25; Value = 827},
    @{Row = 26; Value = 958},
    @{Row = 27; Value = 899},
    @{Row = 28; Value = 964},
    @{Row = 29; Value = 1010},
    @{Row = 30; Value = 1056},
    @{Row = 31; Value = 843},
    @{Row = 32; Value = 801},
    @{Row = 33; Value = 969},
    @{Row = 34; Value = 818},
    @{Row = 35; Value = 942},
    @{Row = 36; Value = 1072},
    @{Row = 37; Value = 982},
    @{Row = 38; Value = 1001},
    @{Row = 39; Value = 1023},
    @{Row = 40; Value = 1988},
    @{Row = 41; Value = 833},
    @{Row = 42; Value = 1024},
    @{Row = 43; Value = 1093},
    @{Row = 44; Value = 969},
    @{Row = 45; Value = 1043},
    @{Row = 46; Value = 932},
    @{Row = 47; Value = 1104},
    @{Row = 48; Value = 872},
    @{Row = 49; Value = 1345},
    @{Row = 50; Value = 903},
    @{Row = 51; Value = 866},
    @{Row = 52; Value = 1221},
    @{Row = 53; Value = 806},
    @{Row = 54; Value = 1884},
    @{Row = 55; Value = 949},
    @{Row = 56; Value = 1015},
    @{Row = 57; Value = 3262}
)

foreach ($u in $supUpdates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.Value
}
